$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2664.1667
$ws.Cells.Item(100, 9).Value = 1715.7
$ws.Cells.Item(100, 10).Value = 3849.75
$ws.Cells.Item(100, 11).Value = 1715.7
$ws.Cells.Item(100, 12).Value = 3849.75
$ws.Cells.Item(100, 13).Value = -1174.7
$ws.Cells.Item(100, 14).Value = -4931.75

$ws.Cells.Item(103, 8).Value = 865.381
$ws.Cells.Item(103, 9).Value = 735.25
$ws.Cells.Item(103, 10).Value = 1038.8889
$ws.Cells.Item(103, 11).Value = 2205.75
$ws.Cells.Item(103, 12).Value = 3116.6667
$ws.Cells.Item(103, 13).Value = -1619.75
$ws.Cells.Item(103, 14).Value = -4288.6667

$ws.Cells.Item(132, 8).Value = 31360.969
$ws.Cells.Item(132, 9).Value = 5728.2383
$ws.Cells.Item(132, 10).Value = 76218.25
$ws.Cells.Item(132, 11).Value = 17184.7149
$ws.Cells.Item(132, 12).Value = 228654.75
$ws.Cells.Item(132, 13).Value = -14654.7149
$ws.Cells.Item(132, 14).Value = -233714.75

$ws.Cells.Item(138, 8).Value = 1350.25
$ws.Cells.Item(138, 9).Value = 614.04443
$ws.Cells.Item(138, 10).Value = 1952.6
$ws.Cells.Item(138, 11).Value = 1842.13329
$ws.Cells.Item(138, 12).Value = 5857.799999999999
$ws.Cells.Item(138, 13).Value = 3297.86671
$ws.Cells.Item(138, 14).Value = -16137.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10439.1455
$ws.Cells.Item(32, 9).Value = 9883.642
$ws.Cells.Item(32, 10).Value = 13710.444
$ws.Cells.Item(32, 11).Value = 9883.642
$ws.Cells.Item(32, 12).Value = 13710.444
$ws.Cells.Item(32, 13).Value = -9596.642
$ws.Cells.Item(32, 14).Value = -14284.444

$ws.Cells.Item(61, 8).Value = 1520.5385
$ws.Cells.Item(61, 9).Value = 1183
$ws.Cells.Item(61, 10).Value = 2158.111
$ws.Cells.Item(61, 11).Value = 1183
$ws.Cells.Item(61, 12).Value = 2158.111
$ws.Cells.Item(61, 13).Value = -971
$ws.Cells.Item(61, 14).Value = -2582.111

$ws.Cells.Item(110, 8).Value = 1487.4048
$ws.Cells.Item(110, 9).Value = 1500.2424
$ws.Cells.Item(110, 10).Value = 1440.3334
$ws.Cells.Item(110, 11).Value = 1500.2424
$ws.Cells.Item(110, 12).Value = 1440.3334
$ws.Cells.Item(110, 13).Value = 544.7575999999999
$ws.Cells.Item(110, 14).Value = -5530.3334

$ws.Cells.Item(115, 8).Value = 20000
$ws.Cells.Item(115, 10).Value = 20000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 14).Value = -23134

$ws.Cells.Item(122, 8).Value = 1769.6666
$ws.Cells.Item(122, 9).Value = 1762.9286
$ws.Cells.Item(122, 10).Value = 1793.25
$ws.Cells.Item(122, 11).Value = 5288.7858
$ws.Cells.Item(122, 12).Value = 5379.75
$ws.Cells.Item(122, 13).Value = -2838.7858
$ws.Cells.Item(122, 14).Value = -10279.75

$ws.Cells.Item(136, 8).Value = 1520.5385
$ws.Cells.Item(136, 9).Value = 1183
$ws.Cells.Item(136, 10).Value = 2158.111
$ws.Cells.Item(136, 11).Value = 3549
$ws.Cells.Item(136, 12).Value = 6474.333
$ws.Cells.Item(136, 13).Value = -999
$ws.Cells.Item(136, 14).Value = -11574.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 668.8
$ws.Cells.Item(16, 9).Value = 691.6667
$ws.Cells.Item(16, 10).Value = 653.55554
$ws.Cells.Item(16, 11).Value = 691.6667
$ws.Cells.Item(16, 12).Value = 653.55554
$ws.Cells.Item(16, 13).Value = -404.6667
$ws.Cells.Item(16, 14).Value = -1227.55554

$ws.Cells.Item(99, 8).Value = 2485.3635
$ws.Cells.Item(99, 9).Value = 1817.5
$ws.Cells.Item(99, 10).Value = 3286.8
$ws.Cells.Item(99, 11).Value = 1817.5
$ws.Cells.Item(99, 12).Value = 3286.8
$ws.Cells.Item(99, 13).Value = -319.5
$ws.Cells.Item(99, 14).Value = -6282.8

$ws.Cells.Item(113, 8).Value = 668.8
$ws.Cells.Item(113, 9).Value = 691.6667
$ws.Cells.Item(113, 10).Value = 653.55554
$ws.Cells.Item(113, 11).Value = 691.6667
$ws.Cells.Item(113, 12).Value = 653.55554
$ws.Cells.Item(113, 13).Value = 1478.3333
$ws.Cells.Item(113, 14).Value = -4993.55554

$ws.Cells.Item(126, 8).Value = 2485.3635
$ws.Cells.Item(126, 9).Value = 1817.5
$ws.Cells.Item(126, 10).Value = 3286.8
$ws.Cells.Item(126, 11).Value = 5452.5
$ws.Cells.Item(126, 12).Value = 9860.400000000001
$ws.Cells.Item(126, 13).Value = -2982.5
$ws.Cells.Item(126, 14).Value = -14800.4

$ws.Cells.Item(132, 8).Value = 144920.7
$ws.Cells.Item(132, 9).Value = 1999
$ws.Cells.Item(132, 10).Value = 180651.12
$ws.Cells.Item(132, 11).Value = 5997
$ws.Cells.Item(132, 12).Value = 541953.36
$ws.Cells.Item(132, 13).Value = -3467
$ws.Cells.Item(132, 14).Value = -547013.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2300.2144
$ws.Cells.Item(34, 10).Value = 2300.2144
$ws.Cells.Item(34, 12).Value = 6900.6432
$ws.Cells.Item(34, 14).Value = -7068.6432

$ws.Cells.Item(98, 8).Value = 1433.3334
$ws.Cells.Item(98, 10).Value = 1433.3334
$ws.Cells.Item(98, 12).Value = 4300.0002
$ws.Cells.Item(98, 14).Value = -7296.0002

$ws.Cells.Item(113, 8).Value = 1972.8572
$ws.Cells.Item(113, 9).Value = 3000.05
$ws.Cells.Item(113, 10).Value = 603.26666
$ws.Cells.Item(113, 11).Value = 9000.150000000001
$ws.Cells.Item(113, 12).Value = 1809.79998
$ws.Cells.Item(113, 13).Value = -6830.150000000001
$ws.Cells.Item(113, 14).Value = -6149.79998

$ws.Cells.Item(122, 8).Value = 7800.067
$ws.Cells.Item(122, 9).Value = 404.66666
$ws.Cells.Item(122, 10).Value = 18893.166
$ws.Cells.Item(122, 11).Value = 3641.99994
$ws.Cells.Item(122, 12).Value = 170038.494
$ws.Cells.Item(122, 13).Value = -1191.99994
$ws.Cells.Item(122, 14).Value = -174938.494

$ws.Cells.Item(131, 8).Value = 5450.5186
$ws.Cells.Item(131, 10).Value = 2414.1667
$ws.Cells.Item(131, 12).Value = 7242.500100000001
$ws.Cells.Item(131, 14).Value = -17322.5001

$ws.Cells.Item(132, 8).Value = 1991.2963
$ws.Cells.Item(132, 10).Value = 2182.7368
$ws.Cells.Item(132, 12).Value = 19644.6312
$ws.Cells.Item(132, 14).Value = -24704.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3175.8
$ws.Cells.Item(80, 9).Value = 3062.4092
$ws.Cells.Item(80, 10).Value = 3367.6924
$ws.Cells.Item(80, 11).Value = 3062.4092
$ws.Cells.Item(80, 12).Value = 3367.6924
$ws.Cells.Item(80, 13).Value = -2064.4092
$ws.Cells.Item(80, 14).Value = -5363.6924

$ws.Cells.Item(83, 8).Value = 3175.8
$ws.Cells.Item(83, 9).Value = 3062.4092
$ws.Cells.Item(83, 10).Value = 3367.6924
$ws.Cells.Item(83, 11).Value = 15312.046
$ws.Cells.Item(83, 12).Value = 16838.462
$ws.Cells.Item(83, 13).Value = -10320.046
$ws.Cells.Item(83, 14).Value = -26822.462

$ws.Cells.Item(123, 8).Value = 17920.223
$ws.Cells.Item(123, 10).Value = 17920.223
$ws.Cells.Item(123, 12).Value = 17920.223
$ws.Cells.Item(123, 14).Value = -22820.223

$ws.Cells.Item(132, 8).Value = 3296.2222
$ws.Cells.Item(132, 9).Value = 2343.5
$ws.Cells.Item(132, 10).Value = 4322.231
$ws.Cells.Item(132, 11).Value = 7030.5
$ws.Cells.Item(132, 12).Value = 12966.693
$ws.Cells.Item(132, 13).Value = -4500.5
$ws.Cells.Item(132, 14).Value = -18026.693

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(121, 8).Value = 40000
$ws.Cells.Item(121, 10).Value = 40000
$ws.Cells.Item(121, 12).Value = 40000
$ws.Cells.Item(121, 14).Value = -43494

$ws.Cells.Item(132, 8).Value = 2092.3247
$ws.Cells.Item(132, 9).Value = 1518.3091
$ws.Cells.Item(132, 10).Value = 3527.3635
$ws.Cells.Item(132, 11).Value = 4554.927299999999
$ws.Cells.Item(132, 12).Value = 10582.0905
$ws.Cells.Item(132, 13).Value = -2024.927299999999
$ws.Cells.Item(132, 14).Value = -15642.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 49142
$ws.Cells.Item(123, 10).Value = 49142
$ws.Cells.Item(123, 12).Value = 49142
$ws.Cells.Item(123, 14).Value = -58942
